# Apply crypto price/volume updates from the Fri Sep 29 20:43:49 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.58"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.667.23"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = "  +4.95%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.249"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.26"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.63%  "
$ws.Range("D12").Value = "1.901.03"
$ws.Range("D13").Value = "1.663.20"
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("D17").Value = "26.913.95"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.94"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.98"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  +0.50%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.115"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.89"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("D33").Value = "1.456.01"
$ws.Range("E33").Value = "  -3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.581"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0170"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "66.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.970"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.46%  "
$ws.Range("D45").Value = "1.809.30"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.52"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.54"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("E50").Value = "  +4.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.00%  "
